$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that used to sit after
#    "What is it?" (it is simply dropped in the target revision).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Locate the "Map/peek" bullet and turn it into the italic
#    "<all RxMobile operators>" line, split across three runs so the
#    "RxMobile" run is bracketed by proofErr spell-check markers, and
#    the whole paragraph gets the extra <w:i/> run-property.
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Map/peek`r") {
        $target = $p
        break
    }
}
if (-not $target) {
    throw "Could not find the 'Map/peek' paragraph"
}

$rng = $target.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p w:rsidR="00994BF8" w:rsidRDefault="00994BF8" w:rsidP="00994BF8">' +
  '<w:pPr>' +
    '<w:pStyle w:val="Lijstalinea"/>' +
    '<w:numPr><w:ilvl w:val="3"/><w:numId w:val="2"/></w:numPr>' +
    '<w:rPr><w:i/><w:lang w:val="en-US"/></w:rPr>' +
  '</w:pPr>' +
  '<w:r><w:rPr><w:i/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">&lt;all </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:i/><w:lang w:val="en-US"/></w:rPr><w:t>RxMobile</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:i/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> operators&gt;</w:t></w:r>' +
'</w:p>' +
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml)

# ---------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark so it spans the freshly inserted
#    paragraph (this is where Word leaves the cursor after the edit).
# ---------------------------------------------------------------------
$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "<all RxMobile operators>`r") {
        $newPara = $p
        break
    }
}
if (-not $newPara) {
    throw "Could not find the freshly inserted paragraph"
}
$d.Bookmarks.Add("_GoBack", $newPara.Range)
